$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3681.0344
$ws.Range("I74").Value = 3573.75
$ws.Range("J74").Value = 4196
$ws.Range("K74").Value = 3573.75
$ws.Range("L74").Value = 4196
$ws.Range("M74").Value = -2637.75
$ws.Range("N74").Value = -6068
$ws.Range("H76").Value = 90911700
$ws.Range("I76").Value = 100002580
$ws.Range("J76").Value = 2900
$ws.Range("K76").Value = 100002580
$ws.Range("L76").Value = 2900
$ws.Range("M76").Value = -100002265
$ws.Range("N76").Value = -3530
$ws.Range("H77").Value = 3681.0344
$ws.Range("I77").Value = 3573.75
$ws.Range("J77").Value = 4196
$ws.Range("K77").Value = 17868.75
$ws.Range("L77").Value = 20980
$ws.Range("M77").Value = -13188.75
$ws.Range("N77").Value = -30340
$ws.Range("H79").Value = 90911700
$ws.Range("I79").Value = 100002580
$ws.Range("J79").Value = 2900
$ws.Range("K79").Value = 100002580
$ws.Range("L79").Value = 2900
$ws.Range("M79").Value = -100001488
$ws.Range("N79").Value = -5084
$ws.Range("H141").Value = 2041.7084
$ws.Range("I141").Value = 1704.3334
$ws.Range("K141").Value = 5113.0002
$ws.Range("M141").Value = 66.9997999999996
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1683.3334
$ws.Range("I63").Value = 1664.7059
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 1664.7059
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = -978.7058999999999
$ws.Range("N63").Value = -3372
$ws.Range("H66").Value = 1683.3334
$ws.Range("I66").Value = 1664.7059
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 8323.529500000001
$ws.Range("L66").Value = 10000
$ws.Range("M66").Value = -4891.529500000001
$ws.Range("N66").Value = -16864
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 492.07693
$ws.Range("I22").Value = 492.07693
$ws.Range("K22").Value = 492.07693
$ws.Range("M22").Value = -319.07693
$ws.Range("H105").Value = 1259.875
$ws.Range("I105").Value = 1279.35
$ws.Range("J105").Value = 1162.5
$ws.Range("K105").Value = 1279.35
$ws.Range("L105").Value = 1162.5
$ws.Range("M105").Value = 467.6500000000001
$ws.Range("N105").Value = -4656.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2609195.2
$ws.Range("I31").Value = 5210798.5
$ws.Range("J31").Value = 7591.8125
$ws.Range("K31").Value = 5210798.5
$ws.Range("L31").Value = 7591.8125
$ws.Range("M31").Value = -5210503.5
$ws.Range("N31").Value = -8181.8125
$ws.Range("H34").Value = 2609195.2
$ws.Range("I34").Value = 5210798.5
$ws.Range("J34").Value = 7591.8125
$ws.Range("K34").Value = 5210798.5
$ws.Range("L34").Value = 7591.8125
$ws.Range("M34").Value = -5210596.5
$ws.Range("N34").Value = -7995.8125
$ws.Range("H50").Value = 13154.375
$ws.Range("J50").Value = 13154.375
$ws.Range("L50").Value = 13154.375
$ws.Range("N50").Value = -14404.375
$ws.Range("H51").Value = 9200.429
$ws.Range("J51").Value = 9200.429
$ws.Range("L51").Value = 9200.429
$ws.Range("N51").Value = -10672.429
$ws.Range("H59").Value = 16500.334
$ws.Range("J59").Value = 16500.334
$ws.Range("L59").Value = 16500.334
$ws.Range("N59").Value = -18790.334
$ws.Range("H60").Value = 4218.4546
$ws.Range("I60").Value = 1942.8572
$ws.Range("J60").Value = 8200.75
$ws.Range("K60").Value = 1942.8572
$ws.Range("L60").Value = 8200.75
$ws.Range("M60").Value = -1431.8572
$ws.Range("N60").Value = -9222.75
$ws.Range("H61").Value = 9200.429
$ws.Range("J61").Value = 9200.429
$ws.Range("L61").Value = 9200.429
$ws.Range("N61").Value = -9896.429
$ws.Range("H62").Value = 2523.3125
$ws.Range("I62").Value = 1817.3
$ws.Range("J62").Value = 3700
$ws.Range("K62").Value = 1817.3
$ws.Range("L62").Value = 3700
$ws.Range("M62").Value = -1193.3
$ws.Range("N62").Value = -4948
$ws.Range("H65").Value = 2523.3125
$ws.Range("I65").Value = 1817.3
$ws.Range("J65").Value = 3700
$ws.Range("K65").Value = 9086.5
$ws.Range("L65").Value = 18500
$ws.Range("M65").Value = -5966.5
$ws.Range("N65").Value = -24740
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 16499.666
$ws.Range("J57").Value = 16499.666
$ws.Range("L57").Value = 16499.666
$ws.Range("N57").Value = -18139.666
$ws.Range("H70").Value = 7085676.5
$ws.Range("I70").Value = 3681070
$ws.Range("J70").Value = 11908870
$ws.Range("K70").Value = 3681070
$ws.Range("L70").Value = 11908870
$ws.Range("M70").Value = -3680800
$ws.Range("N70").Value = -11909410
$ws.Range("H73").Value = 7085676.5
$ws.Range("I73").Value = 3681070
$ws.Range("J73").Value = 11908870
$ws.Range("K73").Value = 3681070
$ws.Range("L73").Value = 11908870
$ws.Range("M73").Value = -3680134
$ws.Range("N73").Value = -11910742
$ws.Range("H80").Value = 9509.462
$ws.Range("I80").Value = 4409.3335
$ws.Range("J80").Value = 16464.182
$ws.Range("K80").Value = 4409.3335
$ws.Range("L80").Value = 16464.182
$ws.Range("M80").Value = -3411.3335
$ws.Range("N80").Value = -18460.182
$ws.Range("H83").Value = 9509.462
$ws.Range("I83").Value = 4409.3335
$ws.Range("J83").Value = 16464.182
$ws.Range("K83").Value = 22046.6675
$ws.Range("L83").Value = 82320.91
$ws.Range("M83").Value = -17054.6675
$ws.Range("N83").Value = -92304.91
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2431.875
$ws.Range("I7").Value = 1475
$ws.Range("J7").Value = 2750.8333
$ws.Range("K7").Value = 1475
$ws.Range("L7").Value = 2750.8333
$ws.Range("M7").Value = -1363
$ws.Range("N7").Value = -2974.8333
$ws.Range("H122").Value = 2610.2083
$ws.Range("I122").Value = 2573.6667
$ws.Range("J122").Value = 2866
$ws.Range("K122").Value = 7721.000100000001
$ws.Range("L122").Value = 8598
$ws.Range("M122").Value = -5271.000100000001
$ws.Range("N122").Value = -13498
$ws.Range("H126").Value = 2431.875
$ws.Range("I126").Value = 1475
$ws.Range("J126").Value = 2750.8333
$ws.Range("K126").Value = 4425
$ws.Range("L126").Value = 8252.499899999999
$ws.Range("M126").Value = -1955
$ws.Range("N126").Value = -13192.4999
$ws.Range("H139").Value = 38533.332
$ws.Range("J139").Value = 38533.332
$ws.Range("L139").Value = 38533.332
$ws.Range("N139").Value = -48813.332
